$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F22").Value = "AUC"
$ws.Range("F23").Value = "AUC"
$ws.Range("D24").Value = 0.9999999999999999
$ws.Range("F24").Value = "AUC"
$ws.Range("F25").Value = "AUC"
$ws.Range("F26").Value = "AUC"
$ws.Range("F27").Value = "AUC"
$ws.Range("A28").Value = 0.9999999999999999
$ws.Range("F28").Value = "AUC"
$ws.Range("C29").Value = 0.9999999999999999
$ws.Range("F29").Value = "AUC"
$ws.Range("F30").Value = "AUC"
$ws.Range("C31").Value = 0.9999999999999999
$ws.Range("F31").Value = "AUC"
$ws.Range("A32").Value = 0.817540847657428
$ws.Range("B32").Value = 0.8669375565424788
$ws.Range("C32").Value = 0.8685841420621214
$ws.Range("D32").Value = 0.819241405543219
$ws.Range("F32").Value = "AUC"
$ws.Range("A33").Value = 0.8133627011541519
$ws.Range("B33").Value = 0.8592597732817941
$ws.Range("C33").Value = 0.8560514159607424
$ws.Range("D33").Value = 0.7988244304630315
$ws.Range("F33").Value = "AUC"
$ws.Range("A34").Value = 0.8038182190169301
$ws.Range("B34").Value = 0.8257471911070192
$ws.Range("C34").Value = 0.823181252450855
$ws.Range("D34").Value = 0.7818466233611239
$ws.Range("F34").Value = "AUC"
$ws.Range("A35").Value = 0.8444887729527901
$ws.Range("B35").Value = 0.8573764342829864
$ws.Range("C35").Value = 0.857986803744054
$ws.Range("D35").Value = 0.8167058803471264
$ws.Range("F35").Value = "AUC"
$ws.Range("A36").Value = 0.8291195676265493
$ws.Range("B36").Value = 0.8477094095783677
$ws.Range("C36").Value = 0.8376928716348694
$ws.Range("D36").Value = 0.8025915127955945
$ws.Range("F36").Value = "AUC"
$ws.Range("A37").Value = 0.8252877090685897
$ws.Range("B37").Value = 0.8528685659727551
$ws.Range("C37").Value = 0.8505421717559205
$ws.Range("D37").Value = 0.7917387004927284
$ws.Range("F37").Value = "AUC"
$ws.Range("A38").Value = 0.8168218164458767
$ws.Range("B38").Value = 0.8651023818048522
$ws.Range("C38").Value = 0.8645568001636745
$ws.Range("D38").Value = 0.7966165413533833
$ws.Range("F38").Value = "AUC"
$ws.Range("A39").Value = 0.8252237737200143
$ws.Range("B39").Value = 0.8790104512983137
$ws.Range("C39").Value = 0.870408163265306
$ws.Range("D39").Value = 0.8421555589654408
$ws.Range("F39").Value = "AUC"
$ws.Range("A40").Value = 0.8666589944248376
$ws.Range("B40").Value = 0.9025096755494176
$ws.Range("C40").Value = 0.9041276661040356
$ws.Range("D40").Value = 0.8331841508533239
$ws.Range("F40").Value = "AUC"
$ws.Range("A41").Value = 0.8581540956014555
$ws.Range("B41").Value = 0.8866101744300643
$ws.Range("C41").Value = 0.8877431042963924
$ws.Range("D41").Value = 0.836809251289329
$ws.Range("F41").Value = "AUC"
$ws.Range("A42").Value = -0.0007834746323563267
$ws.Range("B42").Value = -0.0004371154131972546
$ws.Range("C42").Value = -0.0006668655307575344
$ws.Range("D42").Value = -0.0005400911875440716
$ws.Range("F42").Value = "log loss"
$ws.Range("A43").Value = -0.0003969153974189798
$ws.Range("B43").Value = -0.0005101248051727581
$ws.Range("C43").Value = -0.0005647992604142721
$ws.Range("D43").Value = -0.0004734058803413595
$ws.Range("F43").Value = "log loss"
$ws.Range("A44").Value = -0.0004345655301062264
$ws.Range("B44").Value = -0.0004035605680481028
$ws.Range("C44").Value = -0.0004542955031632838
$ws.Range("D44").Value = -0.0005357143434201267
$ws.Range("F44").Value = "log loss"
$ws.Range("A45").Value = -0.0003590763276736444
$ws.Range("B45").Value = -0.000413723710588281
$ws.Range("C45").Value = -0.000788257174575135
$ws.Range("D45").Value = -0.0003858812574854307
$ws.Range("F45").Value = "log loss"
$ws.Range("A46").Value = -0.0004178049407673643
$ws.Range("B46").Value = -0.0002982008263123635
$ws.Range("C46").Value = -0.0009739946674781173
$ws.Range("D46").Value = -0.0005694326350820679
$ws.Range("F46").Value = "log loss"
$ws.Range("A47").Value = -0.0004703871533440037
$ws.Range("B47").Value = -0.0003488519342311819
$ws.Range("C47").Value = -0.0005928242316756746
$ws.Range("D47").Value = -0.0005294177855590024
$ws.Range("F47").Value = "log loss"
$ws.Range("A48").Value = -0.0004079153759773753
$ws.Range("B48").Value = -0.0004619052853879943
$ws.Range("C48").Value = -0.0009003213654714627
$ws.Range("D48").Value = -0.0004810170234429494
$ws.Range("F48").Value = "log loss"
$ws.Range("A49").Value = -0.0003385158698051249
$ws.Range("B49").Value = -0.0004109182066639176
$ws.Range("C49").Value = -0.0006547624827905523
$ws.Range("D49").Value = -0.0004030655356759185
$ws.Range("F49").Value = "log loss"
$ws.Range("A50").Value = -0.0007331929347144702
$ws.Range("B50").Value = -0.0006867943468819407
$ws.Range("C50").Value = -0.002052529142241044
$ws.Range("D50").Value = -0.0009452307090967786
$ws.Range("F50").Value = "log loss"
$ws.Range("A51").Value = -0.0006567374207898593
$ws.Range("B51").Value = -0.0007775325125011385
$ws.Range("C51").Value = -0.001147144579002779
$ws.Range("D51").Value = -0.0008937360224664204
$ws.Range("F51").Value = "log loss"
$ws.Range("A52").Value = -0.2425772749262943
$ws.Range("B52").Value = -0.1560163855876009
$ws.Range("C52").Value = -0.1551385649568812
$ws.Range("D52").Value = -0.1715481444733716
$ws.Range("F52").Value = "log loss"
$ws.Range("A53").Value = -0.2559245721234976
$ws.Range("B53").Value = -0.1572683465618068
$ws.Range("C53").Value = -0.1576150138725978
$ws.Range("D53").Value = -0.1657465342226435
$ws.Range("F53").Value = "log loss"
$ws.Range("A54").Value = -0.2632588875092604
$ws.Range("B54").Value = -0.1643754880401221
$ws.Range("C54").Value = -0.1631212811610824
$ws.Range("D54").Value = -0.1747094658981536
$ws.Range("F54").Value = "log loss"
$ws.Range("A55").Value = -0.2250558635004151
$ws.Range("B55").Value = -0.1595197546985891
$ws.Range("C55").Value = -0.1576410364645288
$ws.Range("D55").Value = -0.1675316150967985
$ws.Range("F55").Value = "log loss"
$ws.Range("A56").Value = -0.2372810302393573
$ws.Range("B56").Value = -0.1657850097608822
$ws.Range("C56").Value = -0.1685943439614503
$ws.Range("D56").Value = -0.1706282029611253
$ws.Range("F56").Value = "log loss"
$ws.Range("A57").Value = -0.2514967118505977
$ws.Range("B57").Value = -0.1577106252078742
$ws.Range("C57").Value = -0.1573033080164973
$ws.Range("D57").Value = -0.1728537395840315
$ws.Range("F57").Value = "log loss"
$ws.Range("A58").Value = -0.2447620765013136
$ws.Range("B58").Value = -0.1547499616141116
$ws.Range("C58").Value = -0.1542666551808109
$ws.Range("D58").Value = -0.1680923461020287
$ws.Range("F58").Value = "log loss"
$ws.Range("A59").Value = -0.2377780352155519
$ws.Range("B59").Value = -0.1505491647602472
$ws.Range("C59").Value = -0.1556256236738262
$ws.Range("D59").Value = -0.1646216239942101
$ws.Range("F59").Value = "log loss"
$ws.Range("A60").Value = -0.2062961421409688
$ws.Range("B60").Value = -0.1456687701409266
$ws.Range("C60").Value = -0.1449686015754659
$ws.Range("D60").Value = -0.1643484120901253
$ws.Range("F60").Value = "log loss"
$ws.Range("A61").Value = -0.2043559304208107
$ws.Range("B61").Value = -0.1496997405738351
$ws.Range("C61").Value = -0.1496944343718555
$ws.Range("D61").Value = -0.1640823944256502
$ws.Range("F61").Value = "log loss"
